$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 - this shifts existing rows 59..102 down to 60..103
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the new record's data
$ws.Cells.Item(59, 1).Value = 5
$ws.Cells.Item(59, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(59, 3).Value = "Maule"
$ws.Cells.Item(59, 4).Value = 45126
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(59, 6).Value = 100112040
$ws.Cells.Item(59, 7).Value = "Cilantro"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 150
$ws.Cells.Item(59, 11).Value = 12000
$ws.Cells.Item(59, 12).Value = 12000
$ws.Cells.Item(59, 13).Value = 12000
$ws.Cells.Item(59, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(59, 15).Value = "Región Metropolitana"
$ws.Cells.Item(59, 16).Value = 333
$ws.Cells.Item(59, 17).Value = 36
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# Apply the same date number format (s="2") used by the other column D cells
$ws.Cells.Item(59, 4).NumberFormat = $ws.Cells.Item(60, 4).NumberFormat
